$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column D, shifting D:K -> E:L
$ws.Columns.Item(4).Insert()

# Copy cell formatting from the (now-shifted) former column D, now in column E,
# into the freshly inserted column D so the new column matches styling
# (date format row, number format rows, etc.)
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new column D with the latest reporting-period figures.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 64500
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = -700
$ws.Range("D17").Value = 9000
$ws.Range("D18").Value = 55500
$ws.Range("D20").Value = -28800
$ws.Range("D21").Value = 29600
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 26700
$ws.Range("D24").Value = 5000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 21700
$ws.Range("D27").Value = 21700
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 28800
$ws.Range("D33").Value = 21700
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 21700
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 20100
$ws.Range("D42").Value = 24500
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 26400
$ws.Range("D49").Value = 24000
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1647700
$ws.Range("D57").Value = 0
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1479600
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 121900
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 168100
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 21700
$ws.Range("D83").Value = 2900
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 29500
$ws.Range("D91").Value = -1700
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -53600
$ws.Range("D96").Value = -6300
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 30500
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 6500
